# Update "想去人数" (want-to-go count) figures across the workbook.
# Sheets affected: 展览 (Exhibitions), 演出 (Performances), 全部类型 (All types, a merged view).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value  = 8217
$wsExhibit.Range("F5").Value  = 5989
$wsExhibit.Range("F6").Value  = 513
$wsExhibit.Range("F10").Value = 305
$wsExhibit.Range("F11").Value = 840

# --- Sheet: 演出 ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 94

# --- Sheet: 全部类型 (merged list containing rows from both sheets above) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 8217
$wsAll.Range("F5").Value  = 5989
$wsAll.Range("F6").Value  = 513
$wsAll.Range("F10").Value = 305
$wsAll.Range("F11").Value = 94
$wsAll.Range("F15").Value = 840
